$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 25, shifting existing rows 25-103 down to 26-104
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with the new record
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44414
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = 100112032
$ws.Range("G25").Value = "Zapallo italiano"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = 9500
$ws.Range("N25").Value = "$/caja 50 unidades"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 190
$ws.Range("Q25").Value = 50
$ws.Range("R25").Value = "Hortaliza"
